# Updated cryptos list on Thu Dec 14 19:47:21 UTC 2023 with GitHub Actions
#
# Refresh the Price (col D) and Volume(1h) (col E) values for every coin
# row. The source data stores these as plain text, so each value below is
# written with a leading apostrophe -- this forces Excel to keep it as
# text instead of auto-converting decimal-looking strings (e.g. "235.70",
# "3.90") into numbers, which would silently drop the trailing zero.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.902.83"
$ws.Range("E2").Value = "'  +0.18%  "
$ws.Range("D3").Value = "'2.288.75"
$ws.Range("E3").Value = "'  +1.59%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'251.87"
$ws.Range("E5").Value = "'  -0.54%  "
$ws.Range("D6").Value = "'0.641"
$ws.Range("E6").Value = "'  +0.95%  "
$ws.Range("D7").Value = "'73.75"
$ws.Range("E7").Value = "'  +4.39%  "
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("D9").Value = "'0.647"
$ws.Range("E9").Value = "'  +0.87%  "
$ws.Range("D10").Value = "'39.12"
$ws.Range("E10").Value = "'  -5.16%  "
$ws.Range("D11").Value = "'0.0971"
$ws.Range("E11").Value = "'  +1.18%  "
$ws.Range("D12").Value = "'59.03"
$ws.Range("E12").Value = "'  -0.89%  "
$ws.Range("E13").Value = "'  +1.37%  "
$ws.Range("E14").Value = "'  +0.78%  "
$ws.Range("D15").Value = "'2.631.01"
$ws.Range("E15").Value = "'  +1.72%  "
$ws.Range("D16").Value = "'15.36"
$ws.Range("E16").Value = "'  +3.57%  "
$ws.Range("D17").Value = "'0.871"
$ws.Range("E17").Value = "'  -2.21%  "
$ws.Range("D18").Value = "'2.279.82"
$ws.Range("E18").Value = "'  +1.79%  "
$ws.Range("D19").Value = "'42.796.45"
$ws.Range("E19").Value = "'  +0.20%  "
$ws.Range("E20").Value = "'  +3.00%  "
$ws.Range("E21").Value = "'  +0.94%  "
$ws.Range("D22").Value = "'72.62"
$ws.Range("E22").Value = "'  -0.57%  "
$ws.Range("D23").Value = "'235.70"
$ws.Range("E23").Value = "'  +0.18%  "
$ws.Range("E24").Value = "'  +5.70%  "
$ws.Range("D25").Value = "'3.90"
$ws.Range("E25").Value = "'  -2.32%  "
$ws.Range("D26").Value = "'11.62"
$ws.Range("E26").Value = "'  +0.02%  "
$ws.Range("E27").Value = "'  -0.38%  "
$ws.Range("D28").Value = "'2.41"
$ws.Range("E28").Value = "'  -1.21%  "
$ws.Range("D29").Value = "'3.64"
$ws.Range("E29").Value = "'  -1.07%  "
$ws.Range("D30").Value = "'2.14"
$ws.Range("E30").Value = "'  -2.95%  "
$ws.Range("E31").Value = "'  -0.40%  "
$ws.Range("E32").Value = "'  +0.46%  "
$ws.Range("D33").Value = "'6.44"
$ws.Range("E33").Value = "'  +5.71%  "
$ws.Range("E34").Value = "'  +3.08%  "
$ws.Range("D35").Value = "'0.0828"
$ws.Range("E35").Value = "'  +4.78%  "
$ws.Range("D36").Value = "'30.96"
$ws.Range("E36").Value = "'  +10.83%  "
$ws.Range("E37").Value = "'  +2.11%  "
$ws.Range("D38").Value = "'4.61"
$ws.Range("E38").Value = "'  +10.25%  "
$ws.Range("E39").Value = "'  +1.70%  "
$ws.Range("E40").Value = "'  -3.31%  "
$ws.Range("D41").Value = "'14.32"
$ws.Range("E41").Value = "'  +13.37%  "
$ws.Range("E42").Value = "'  +2.97%  "
$ws.Range("D43").Value = "'5.93"
$ws.Range("E43").Value = "'  +2.25%  "
$ws.Range("D44").Value = "'0.216"
$ws.Range("E44").Value = "'  +6.41%  "
$ws.Range("D45").Value = "'9.19"
$ws.Range("E45").Value = "'  +4.67%  "
$ws.Range("D46").Value = "'61.82"
$ws.Range("E46").Value = "'  -3.49%  "
$ws.Range("D47").Value = "'4.89"
$ws.Range("E47").Value = "'  -2.67%  "
$ws.Range("D48").Value = "'0.103"
$ws.Range("E48").Value = "'  +1.18%  "
$ws.Range("E49").Value = "'  -0.11%  "
$ws.Range("E50").Value = "'  -1.76%  "
$ws.Range("D51").Value = "'100.29"
$ws.Range("E51").Value = "'  +5.74%  "
